$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'big and tall running pants for men'
$ws.Range("A2").Value = 'girls compression knee sleeve'
$ws.Range("A3").Value = 'football pants with pads'
$ws.Range("A4").Value = 'compression padded shorts'
$ws.Range("A5").Value = 'spandex for men pants'
$ws.Range("A6").Value = 'elastic knee pad'
$ws.Range("A7").Value = 'knee pads gel construction'
$ws.Range("A8").Value = 'capri shorts for men'
$ws.Range("A9").Value = 'football padded pants'
$ws.Range("A10").Value = 'girls volleyball kneepads'
$ws.Range("A11").Value = 'black capris men'
$ws.Range("A12").Value = 'womens basketball pants'
$ws.Range("A13").Value = 'baseball material'
$ws.Range("A14").Value = 'youth sports leggings'
$ws.Range("A15").Value = 'protect knee pads'
$ws.Range("A16").Value = 'knees pad'
$ws.Range("A17").Value = 'black youth knee pads'
$ws.Range("A18").Value = 'exercise pads for knees'
$ws.Range("A19").Value = 'gym shorts above knee for men'
$ws.Range("A20").Value = 'knee sleeve basketball youth'
$ws.Range("A21").Value = 'knee pads construction'
$ws.Range("A22").Value = 'work knee pad'
$ws.Range("A23").Value = 'bjj knee'
$ws.Range("A24").Value = 'knee sleeve bjj'
$ws.Range("A25").Value = 'knee pads under'
$ws.Range("A26").Value = 'mens long cycling pants'
$ws.Range("A27").Value = 'mens pad'
$ws.Range("A28").Value = 'boys long baseball pants'
$ws.Range("A29").Value = 'mens basketball gear'
$ws.Range("A30").Value = 'girl sliding shorts'
$ws.Range("A31").Value = 'calf silicone pads'
$ws.Range("A32").Value = 'compression knee sleeve men pair'
$ws.Range("A33").Value = 'girls black softball pants'
$ws.Range("A34").Value = 'hip protector pads'
$ws.Range("A35").Value = 'flexible work pants for men'
$ws.Range("A36").Value = 'knee protection pads'
$ws.Range("A37").Value = 'lightweight work pants for men'
$ws.Range("A38").Value = 'youth baseball pants long'
$ws.Range("A39").Value = 'knees pads work'
$ws.Range("A40").Value = 'tights boys'
$ws.Range("A41").Value = 'mens shorts long below knee'
$ws.Range("A42").Value = 'knee sleeve youth'
$ws.Range("A43").Value = 'snowboarding pants men'
$ws.Range("A44").Value = 'baseball shorts for men'
$ws.Range("A45").Value = 'baseball mens pants'
$ws.Range("A46").Value = 'knee compression sleeve - reduce strain & swelling'
$ws.Range("A47").Value = 'pads men'
$ws.Range("A48").Value = 'basketball sleeve youth leg'
$ws.Range("A49").Value = 'thigh pads football'
$ws.Range("A50").Value = 'compression volleyball'
$ws.Range("A51").Value = 'leggings for mens'
$ws.Range("A52").Value = 'mens yoga pants'
$ws.Range("A53").Value = 'padded football pants'
$ws.Range("A54").Value = 'spandex capris'
$ws.Range("A55").Value = 'water knee hockey'
$ws.Range("A56").Value = 'compression pants sleeves'
$ws.Range("A57").Value = 'knee sleeve padded'
$ws.Range("A58").Value = 'knees pads for construction'
$ws.Range("A59").Value = 'tight capri'
$ws.Range("A60").Value = 'mens baseball compression shorts'
$ws.Range("A61").Value = 'mens running knee compression'
$ws.Range("A62").Value = 'black football leggings'
$ws.Range("A63").Value = 'knee sleeves basketball youth'
$ws.Range("A64").Value = 'cycling knee pads'
$ws.Range("A65").Value = 'construction knee pad'
$ws.Range("A66").Value = 'compression calf leggings'
$ws.Range("A67").Value = 'baseball youth compression sleeve'
$ws.Range("A68").Value = '6 pairs of leggings'
$ws.Range("A69").Value = 'basketball shorts for men pack of 5'
$ws.Range("A70").Value = 'compression pants youth boys'
$ws.Range("A71").Value = 'impact shorts men'
$ws.Range("A72").Value = 'large knee pad'
$ws.Range("A73").Value = 'baseball compression sleeve'
$ws.Range("A74").Value = 'boys sports leggings'
$ws.Range("A75").Value = 'volleyball spandex pack'
$ws.Range("A76").Value = 'baseball pants youth large'
$ws.Range("A77").Value = 'boys paintball pants'
$ws.Range("A78").Value = 'yoga pants mens'
$ws.Range("A79").Value = 'calf tear compression sleeve'
$ws.Range("A80").Value = 'compression shorts men long length'
$ws.Range("A81").Value = 'yoga hand pads'
$ws.Range("A82").Value = 'knee sleeves with padding'
$ws.Range("A83").Value = 'athletic capri leggings'
$ws.Range("A84").Value = 'pants compression men'
$ws.Range("A85").Value = 'basketball padding'
$ws.Range("A86").Value = 'knee pads for men floor work'
$ws.Range("A87").Value = 'youth knee sleeve wrestling'
$ws.Range("A88").Value = 'professional construction knee pads'
$ws.Range("A89").Value = 'youth basketball'
$ws.Range("A90").Value = 'basketball compression knee sleeve'
$ws.Range("A91").Value = 'black softball pants youth girls'
$ws.Range("A92").Value = 'hex gear wash'
$ws.Range("A93").Value = 'knee construction pads'
$ws.Range("A94").Value = 'youth girls softball pants'
$ws.Range("A95").Value = 'compression tight pants'
$ws.Range("A96").Value = 'male workout leggings'
$ws.Range("A97").Value = 'boys xl baseball pants'
$ws.Range("A98").Value = 'thick leggings for men'
$ws.Range("A99").Value = 'knee pads for'
$ws.Range("A100").Value = 'adult football girdle'
